$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.733.31"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.318.21"
$ws.Range("E3").Value = "  +3.21%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'268.97"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "'94.34"
$ws.Range("E6").Value = "  +7.53%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.621"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").Value = "'44.80"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").Value = "'0.0943"
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").Value = "'8.16"
$ws.Range("E12").Value = "  +6.20%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "2.665.22"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").Value = "'15.51"
$ws.Range("E15").Value = "  +3.66%  "
$ws.Range("D16").Value = "'0.865"
$ws.Range("E16").Value = "  +8.42%  "
$ws.Range("D17").Value = "2.322.72"
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").Value = "43.719.01"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  +3.31%  "
$ws.Range("D20").Value = "'6.45"
$ws.Range("E20").Value = "  +7.34%  "
$ws.Range("D21").Value = "'71.50"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").Value = "'238.24"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("E23").Value = "  -4.07%  "
$ws.Range("E24").Value = "  +10.06%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'11.37"
$ws.Range("E26").Value = "  +4.62%  "
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").Value = "'38.55"
$ws.Range("E30").Value = "  -3.09%  "
$ws.Range("D31").Value = "'22.42"
$ws.Range("E31").Value = "  +7.37%  "
$ws.Range("D32").Value = "'171.99"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").Value = "'0.0899"
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("D34").Value = "'5.54"
$ws.Range("E34").Value = "  +2.17%  "
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("D36").Value = "'0.109"
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").Value = "'3.42"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.34"
$ws.Range("E40").Value = "  +6.15%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.235"
$ws.Range("E41").Value = "  +13.99%  "
$ws.Range("E42").Value = "  +20.80%  "
$ws.Range("D43").Value = "'12.13"
$ws.Range("E43").Value = "  -3.33%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'9.12"
$ws.Range("E45").Value = "  +7.20%  "
$ws.Range("B46").Value = "MultiversX"
$ws.Range("C46").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D46").Value = "'61.89"
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("E47").Value = "  +3.16%  "
$ws.Range("D48").Value = "'100.78"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").Value = "2.545.94"
$ws.Range("E50").Value = "  +2.86%  "
$ws.Range("D51").Value = "'0.420"
$ws.Range("E51").Value = "  -2.09%  "
